# edit.ps1 - applies the "Added Goals Sequence Diagram..." commit changes:
#   1) Widen + retitle the "Rectangle 1" label inside the slide's shape
#      group: "Goals Package" -> "Logic, Command, Goals Package"
#      (width grows from 1542217 EMU to 3230180 EMU; position/height unchanged).
#   2) Re-cache the "11/8/19" -> "11/9/19" datetimeFigureOut placeholder text
#      that is stamped on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Recursively look for a shape named "Rectangle 1" (descending into any
# groups) and retitle / widen it once found.
# ---------------------------------------------------------------------
function Fix-GoalsPackageLabel($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -eq "Rectangle 1") {
            $shp.TextFrame.TextRange.Text = "Logic, Command, Goals Package"
            # 3230180 EMU -> points (1 pt = 12700 EMU). Width is a
            # single-precision float property, so round to the value
            # that round-trips to exactly that EMU count.
            $shp.Width = 254.3449
        }
        if ($shp.Type -eq 6) {
            # msoGroup - recurse into it
            Fix-GoalsPackageLabel $shp.GroupItems
        }
    }
}

$slide = $p.Slides.Item(1)
Fix-GoalsPackageLabel $slide.Shapes

# ---------------------------------------------------------------------
# Recursively re-stamp every "Date Placeholder*" shape's cached
# datetimeFigureOut text from 11/8/19 to 11/9/19.
# ---------------------------------------------------------------------
function Update-DateStamp($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "11/8/19") {
                    $shp.TextFrame.TextRange.Text = "11/9/19"
                }
            }
        }
        if ($shp.Type -eq 6) {
            Update-DateStamp $shp.GroupItems
        }
    }
}

$master = $p.SlideMaster
Update-DateStamp $master.Shapes

$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateStamp $layouts.Item($j).Shapes
}
